# Update the Bmp6-Acvr2a LR-pair sheet with refreshed TPM-derived values.
# For each data row (2-19), refresh the ligand/receptor average & total
# expression values, their derived specificity scores, and the
# downstream edge-expression weights/specificities that are computed
# from them (Q = G*M, R = H*N, S = Q / sum(Q), T = R / sum(R)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 7.739652666666667
$ws.Cells.Item(2, 8).Value = 23.218958
$ws.Cells.Item(2, 9).Value = 0.6488398532974882
$ws.Cells.Item(2, 10).Value = 0.6488398532974882
$ws.Cells.Item(2, 13).Value = 11.651608
$ws.Cells.Item(2, 14).Value = 34.954824
$ws.Cells.Item(2, 15).Value = 0.1892813629236475
$ws.Cells.Item(2, 16).Value = 0.1892813629236474
$ws.Cells.Item(2, 17).Value = 90.17939892815468
$ws.Cells.Item(2, 18).Value = 811.614590353392
$ws.Cells.Item(2, 19).Value = 0.122813291751328
$ws.Cells.Item(2, 20).Value = 0.122813291751328

# Row 3
$ws.Cells.Item(3, 7).Value = 7.739652666666667
$ws.Cells.Item(3, 8).Value = 23.218958
$ws.Cells.Item(3, 9).Value = 0.6488398532974882
$ws.Cells.Item(3, 10).Value = 0.6488398532974882
$ws.Cells.Item(3, 15).Value = 0.4419371310876561
$ws.Cells.Item(3, 16).Value = 0.4419371310876561
$ws.Cells.Item(3, 17).Value = 210.5522922591916
$ws.Cells.Item(3, 18).Value = 1894.970630332724
$ws.Cells.Item(3, 19).Value = 0.2867464233016276
$ws.Cells.Item(3, 20).Value = 0.2867464233016276

# Row 4
$ws.Cells.Item(4, 7).Value = 7.739652666666667
$ws.Cells.Item(4, 8).Value = 23.218958
$ws.Cells.Item(4, 9).Value = 0.6488398532974882
$ws.Cells.Item(4, 10).Value = 0.6488398532974882
$ws.Cells.Item(4, 13).Value = 8.657179
$ws.Cells.Item(4, 14).Value = 25.971537
$ws.Cells.Item(4, 15).Value = 0.1406366091439035
$ws.Cells.Item(4, 16).Value = 0.1406366091439035
$ws.Cells.Item(4, 17).Value = 67.00355853316067
$ws.Cells.Item(4, 18).Value = 603.0320267984459
$ws.Cells.Item(4, 19).Value = 0.09125063684518653
$ws.Cells.Item(4, 20).Value = 0.09125063684518653

# Row 5
$ws.Cells.Item(5, 7).Value = 7.739652666666667
$ws.Cells.Item(5, 8).Value = 23.218958
$ws.Cells.Item(5, 9).Value = 0.6488398532974882
$ws.Cells.Item(5, 10).Value = 0.6488398532974882
$ws.Cells.Item(5, 13).Value = 5.488499666666667
$ws.Cells.Item(5, 14).Value = 16.465499
$ws.Cells.Item(5, 15).Value = 0.08916114387925267
$ws.Cells.Item(5, 16).Value = 0.08916114387925267
$ws.Cells.Item(5, 17).Value = 42.47908108111578
$ws.Cells.Item(5, 18).Value = 382.311729730042
$ws.Cells.Item(5, 19).Value = 0.05785130351445054
$ws.Cells.Item(5, 20).Value = 0.05785130351445054

# Row 6
$ws.Cells.Item(6, 7).Value = 7.739652666666667
$ws.Cells.Item(6, 8).Value = 23.218958
$ws.Cells.Item(6, 9).Value = 0.6488398532974882
$ws.Cells.Item(6, 10).Value = 0.6488398532974882
$ws.Cells.Item(6, 13).Value = 4.091608333333333
$ws.Cells.Item(6, 14).Value = 12.274825
$ws.Cells.Item(6, 15).Value = 0.06646852536431769
$ws.Cells.Item(6, 16).Value = 0.06646852536431769
$ws.Cells.Item(6, 17).Value = 31.66762734803889
$ws.Cells.Item(6, 18).Value = 285.00864613235
$ws.Cells.Item(6, 19).Value = 0.04312742824628426
$ws.Cells.Item(6, 20).Value = 0.04312742824628426

# Row 7
$ws.Cells.Item(7, 7).Value = 7.739652666666667
$ws.Cells.Item(7, 8).Value = 23.218958
$ws.Cells.Item(7, 9).Value = 0.6488398532974882
$ws.Cells.Item(7, 10).Value = 0.6488398532974882
$ws.Cells.Item(7, 13).Value = 4.463825666666667
$ws.Cells.Item(7, 14).Value = 13.391477
$ws.Cells.Item(7, 15).Value = 0.07251522760122259
$ws.Cells.Item(7, 16).Value = 0.07251522760122257
$ws.Cells.Item(7, 17).Value = 34.54846022455178
$ws.Cells.Item(7, 18).Value = 310.936142020966
$ws.Cells.Item(7, 19).Value = 0.04705076963861123
$ws.Cells.Item(7, 20).Value = 0.04705076963861122

# Row 8
$ws.Cells.Item(8, 9).Value = 0.3053032463428815
$ws.Cells.Item(8, 10).Value = 0.3053032463428815
$ws.Cells.Item(8, 13).Value = 11.651608
$ws.Cells.Item(8, 14).Value = 34.954824
$ws.Cells.Item(8, 15).Value = 0.1892813629236475
$ws.Cells.Item(8, 16).Value = 0.1892813629236474
$ws.Cells.Item(8, 17).Value = 42.43275610475201
$ws.Cells.Item(8, 18).Value = 381.8948049427681
$ws.Cells.Item(8, 19).Value = 0.0577882145727947
$ws.Cells.Item(8, 20).Value = 0.05778821457279469

# Row 9
$ws.Cells.Item(9, 9).Value = 0.3053032463428815
$ws.Cells.Item(9, 10).Value = 0.3053032463428815
$ws.Cells.Item(9, 15).Value = 0.4419371310876561
$ws.Cells.Item(9, 16).Value = 0.4419371310876561
$ws.Cells.Item(9, 19).Value = 0.134924840800521
$ws.Cells.Item(9, 20).Value = 0.134924840800521

# Row 10
$ws.Cells.Item(10, 9).Value = 0.3053032463428815
$ws.Cells.Item(10, 10).Value = 0.3053032463428815
$ws.Cells.Item(10, 13).Value = 8.657179
$ws.Cells.Item(10, 14).Value = 25.971537
$ws.Cells.Item(10, 15).Value = 0.1406366091439035
$ws.Cells.Item(10, 16).Value = 0.1406366091439035
$ws.Cells.Item(10, 17).Value = 31.527662539126
$ws.Cells.Item(10, 18).Value = 283.748962852134
$ws.Cells.Item(10, 19).Value = 0.04293681332628871
$ws.Cells.Item(10, 20).Value = 0.04293681332628871

# Row 11
$ws.Cells.Item(11, 9).Value = 0.3053032463428815
$ws.Cells.Item(11, 10).Value = 0.3053032463428815
$ws.Cells.Item(11, 13).Value = 5.488499666666667
$ws.Cells.Item(11, 14).Value = 16.465499
$ws.Cells.Item(11, 15).Value = 0.08916114387925267
$ws.Cells.Item(11, 16).Value = 0.08916114387925267
$ws.Cells.Item(11, 17).Value = 19.98798515506867
$ws.Cells.Item(11, 18).Value = 179.891866395618
$ws.Cells.Item(11, 19).Value = 0.02722118667398058
$ws.Cells.Item(11, 20).Value = 0.02722118667398058

# Row 12
$ws.Cells.Item(12, 9).Value = 0.3053032463428815
$ws.Cells.Item(12, 10).Value = 0.3053032463428815
$ws.Cells.Item(12, 13).Value = 4.091608333333333
$ws.Cells.Item(12, 14).Value = 12.274825
$ws.Cells.Item(12, 15).Value = 0.06646852536431769
$ws.Cells.Item(12, 16).Value = 0.06646852536431769
$ws.Cells.Item(12, 17).Value = 14.90079467868333
$ws.Cells.Item(12, 18).Value = 134.10715210815
$ws.Cells.Item(12, 19).Value = 0.02029305657335035
$ws.Cells.Item(12, 20).Value = 0.02029305657335035

# Row 13
$ws.Cells.Item(13, 9).Value = 0.3053032463428815
$ws.Cells.Item(13, 10).Value = 0.3053032463428815
$ws.Cells.Item(13, 13).Value = 4.463825666666667
$ws.Cells.Item(13, 14).Value = 13.391477
$ws.Cells.Item(13, 15).Value = 0.07251522760122259
$ws.Cells.Item(13, 16).Value = 0.07251522760122257
$ws.Cells.Item(13, 17).Value = 16.25633352991267
$ws.Cells.Item(13, 18).Value = 146.307001769214
$ws.Cells.Item(13, 19).Value = 0.02213913439594618
$ws.Cells.Item(13, 20).Value = 0.02213913439594617

# Row 14
$ws.Cells.Item(14, 7).Value = 0.5470016666666667
$ws.Cells.Item(14, 8).Value = 1.641005
$ws.Cells.Item(14, 9).Value = 0.04585690035963046
$ws.Cells.Item(14, 10).Value = 0.04585690035963046
$ws.Cells.Item(14, 13).Value = 11.651608
$ws.Cells.Item(14, 14).Value = 34.954824
$ws.Cells.Item(14, 15).Value = 0.1892813629236475
$ws.Cells.Item(14, 16).Value = 0.1892813629236474
$ws.Cells.Item(14, 17).Value = 6.373448995346668
$ws.Cells.Item(14, 18).Value = 57.36104095812001
$ws.Cells.Item(14, 19).Value = 0.008679856599524754
$ws.Cells.Item(14, 20).Value = 0.00867985659952475

# Row 15
$ws.Cells.Item(15, 7).Value = 0.5470016666666667
$ws.Cells.Item(15, 8).Value = 1.641005
$ws.Cells.Item(15, 9).Value = 0.04585690035963046
$ws.Cells.Item(15, 10).Value = 0.04585690035963046
$ws.Cells.Item(15, 15).Value = 0.4419371310876561
$ws.Cells.Item(15, 16).Value = 0.4419371310876561
$ws.Cells.Item(15, 17).Value = 14.88082989593222
$ws.Cells.Item(15, 18).Value = 133.92746906339
$ws.Cells.Item(15, 19).Value = 0.02026586698550759
$ws.Cells.Item(15, 20).Value = 0.02026586698550759

# Row 16
$ws.Cells.Item(16, 7).Value = 0.5470016666666667
$ws.Cells.Item(16, 8).Value = 1.641005
$ws.Cells.Item(16, 9).Value = 0.04585690035963046
$ws.Cells.Item(16, 10).Value = 0.04585690035963046
$ws.Cells.Item(16, 13).Value = 8.657179
$ws.Cells.Item(16, 14).Value = 25.971537
$ws.Cells.Item(16, 15).Value = 0.1406366091439035
$ws.Cells.Item(16, 16).Value = 0.1406366091439035
$ws.Cells.Item(16, 17).Value = 4.735491341631667
$ws.Cells.Item(16, 18).Value = 42.619422074685
$ws.Cells.Item(16, 19).Value = 0.006449158972428277
$ws.Cells.Item(16, 20).Value = 0.006449158972428276

# Row 17
$ws.Cells.Item(17, 7).Value = 0.5470016666666667
$ws.Cells.Item(17, 8).Value = 1.641005
$ws.Cells.Item(17, 9).Value = 0.04585690035963046
$ws.Cells.Item(17, 10).Value = 0.04585690035963046
$ws.Cells.Item(17, 13).Value = 5.488499666666667
$ws.Cells.Item(17, 14).Value = 16.465499
$ws.Cells.Item(17, 15).Value = 0.08916114387925267
$ws.Cells.Item(17, 16).Value = 0.08916114387925267
$ws.Cells.Item(17, 17).Value = 3.002218465166111
$ws.Cells.Item(17, 18).Value = 27.019966186495
$ws.Cells.Item(17, 19).Value = 0.004088653690821565
$ws.Cells.Item(17, 20).Value = 0.004088653690821564

# Row 18
$ws.Cells.Item(18, 7).Value = 0.5470016666666667
$ws.Cells.Item(18, 8).Value = 1.641005
$ws.Cells.Item(18, 9).Value = 0.04585690035963046
$ws.Cells.Item(18, 10).Value = 0.04585690035963046
$ws.Cells.Item(18, 13).Value = 4.091608333333333
$ws.Cells.Item(18, 14).Value = 12.274825
$ws.Cells.Item(18, 15).Value = 0.06646852536431769
$ws.Cells.Item(18, 16).Value = 0.06646852536431769
$ws.Cells.Item(18, 17).Value = 2.238116577680556
$ws.Cells.Item(18, 18).Value = 20.143049199125
$ws.Cells.Item(18, 19).Value = 0.003048040544683087
$ws.Cells.Item(18, 20).Value = 0.003048040544683086

# Row 19
$ws.Cells.Item(19, 7).Value = 0.5470016666666667
$ws.Cells.Item(19, 8).Value = 1.641005
$ws.Cells.Item(19, 9).Value = 0.04585690035963046
$ws.Cells.Item(19, 10).Value = 0.04585690035963046
$ws.Cells.Item(19, 13).Value = 4.463825666666667
$ws.Cells.Item(19, 14).Value = 13.391477
$ws.Cells.Item(19, 15).Value = 0.07251522760122259
$ws.Cells.Item(19, 16).Value = 0.07251522760122257
$ws.Cells.Item(19, 17).Value = 2.441720079376112
$ws.Cells.Item(19, 18).Value = 21.975480714385
$ws.Cells.Item(19, 19).Value = 0.003325323566665189
$ws.Cells.Item(19, 20).Value = 0.003325323566665188

